$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column cells we touch to stay as Text so strings like
# "0.9990" / "5.810" / "0.01910" keep their exact digits instead of being
# auto-coerced to numbers (which would strip trailing zeros).

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '29.801.48'
$ws.Range('E2').Value = '  -0.50%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.890.84'
$ws.Range('E3').Value = '  -0.04%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9989'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.7893'
$ws.Range('E5').Value = '  -4.95%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '242.59'
$ws.Range('E6').Value = '  +0.49%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.9992'
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3171'
$ws.Range('E8').Value = '  -1.40%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '25.45'
$ws.Range('E9').Value = '  -4.27%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.07026'
$ws.Range('E10').Value = '  +0.19%  '
$ws.Range('E11').Value = '  +0.12%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.7662'
$ws.Range('E12').Value = '  +2.63%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.884.07'
$ws.Range('E13').Value = '  -0.57%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.287'
$ws.Range('E14').Value = '  +1.85%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '92.15'
$ws.Range('E15').Value = '  -0.11%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '29.799.68'
$ws.Range('E16').Value = '  -0.52%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '13.86'
$ws.Range('E17').Value = '  -1.20%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '5.917'
$ws.Range('E18').Value = '  +0.05%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '243.78'
$ws.Range('E19').Value = '  +0.25%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.000007707'
$ws.Range('E20').Value = '  -0.49%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '8.212'
$ws.Range('E21').Value = '  +18.71%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.9995'
$ws.Range('E22').Value = '  -0.20%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '2.138.63'
$ws.Range('E23').Value = '  -0.74%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.9990'
$ws.Range('E24').Value = '  -0.13%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.1663'
$ws.Range('E25').Value = '  +4.89%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.305'
$ws.Range('E26').Value = '  +1.28%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '165.88'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '18.69'
$ws.Range('E28').Value = '  -0.72%  '
$ws.Range('E29').Value = '  -1.90%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.396'
$ws.Range('E30').Value = '  +1.61%  '
$ws.Range('E31').Value = '  +1.35%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.422'
$ws.Range('E32').Value = '  +4.02%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.05637'
$ws.Range('E33').Value = '  -0.22%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.041'
$ws.Range('E34').Value = '  -0.65%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.261'
$ws.Range('E35').Value = '  -0.97%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.7374'
$ws.Range('E36').Value = '  +0.94%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.9998'
$ws.Range('E37').Value = '  -0.03%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.641'
$ws.Range('E38').Value = '  -2.91%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.01910'
$ws.Range('E39').Value = '  +0.24%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.773'
$ws.Range('E40').Value = '  +0.07%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.4406'
$ws.Range('E41').Value = '  +0.01%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '72.41'
$ws.Range('E42').Value = '  +0.83%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.810'
$ws.Range('E43').Value = '  -2.24%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.025.34'
$ws.Range('E46').Value = '  +3.74%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '102.39'
$ws.Range('E47').Value = '  +1.11%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.866'
$ws.Range('E48').Value = '  -1.11%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '9.901'
$ws.Range('E49').Value = '  +2.39%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.420'
$ws.Range('E50').Value = '  -2.26%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.035.21'
$ws.Range('E51').Value = '  -0.74%  '

# Row 44: now TrustWalletToken (was PaxDollar)
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.8397'
$ws.Range('E44').Value = '  -0.60%  '

# Row 45: now PaxDollar (was TrustWalletToken)
$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.9984'
$ws.Range('E45').Value = '  -0.25%  '
